$wb = $excel.ActiveWorkbook

# Rename sheet "Wong3" to "Euclid"
$ws = $wb.Worksheets.Item("Wong3")
$ws.Name = "Euclid"

# Update SBFL:RANK (column C) and SBFL:EXAM (column D) values
$updates = @(
    @{Row=2;  C=48;  D=19.12350597609562},
    @{Row=3;  C=37;  D=14.74103585657371},
    @{Row=8;  C=9;   D=3.585657370517929},
    @{Row=10; C=37;  D=14.74103585657371},
    @{Row=11; C=8;   D=3.187250996015936},
    @{Row=14; C=165; D=65.73705179282868},
    @{Row=20; C=48;  D=19.12350597609562},
    @{Row=26; C=5;   D=1.99203187250996},
    @{Row=27; C=109; D=43.42629482071713},
    @{Row=28; C=41;  D=16.33466135458167},
    @{Row=29; C=13;  D=5.179282868525896},
    @{Row=33; C=22;  D=8.764940239043826},
    @{Row=35; C=7;   D=2.788844621513944},
    @{Row=36; C=27;  D=10.75697211155378},
    @{Row=40; C=50;  D=19.9203187250996},
    @{Row=42; C=48;  D=19.12350597609562},
    @{Row=45; C=35;  D=13.94422310756972},
    @{Row=48; C=35;  D=13.94422310756972},
    @{Row=49; C=21;  D=8.366533864541832}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
